# Updates cryptos list figures (prices / 1h volume %) per upstream data refresh,
# and fixes the PEPE / Binance-PegBSC-USD row ordering (rows 29-30 swapped).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a literal text value into a cell without Excel re-typing
# numeric-looking strings (e.g. "1.00") as a number.
function Set-TextCell($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 2
$ws.Range("D2").Value = "64.959.77"
$ws.Range("E2").Value = "  +2.99%  "

# Row 3
$ws.Range("D3").Value = "2.629.30"
$ws.Range("E3").Value = "  +2.14%  "

# Row 4
Set-TextCell $ws.Range("D4") "1.00"
$ws.Range("E4").Value = "  +0.02%  "

# Row 5
Set-TextCell $ws.Range("D5") "595.56"
$ws.Range("E5").Value = "  +1.78%  "

# Row 6
Set-TextCell $ws.Range("D6") "155.14"
$ws.Range("E6").Value = "  +5.07%  "

# Row 7
Set-TextCell $ws.Range("D7") "1.00"
$ws.Range("E7").Value = "  +0.01%  "

# Row 8
$ws.Range("E8").Value = "  +1.03%  "

# Row 9
$ws.Range("E9").Value = "  +8.43%  "

# Row 10
$ws.Range("E10").Value = "  +5.27%  "

# Row 11
Set-TextCell $ws.Range("D11") "5.78"
$ws.Range("E11").Value = "  +1.12%  "

# Row 12
$ws.Range("E12").Value = "  +2.18%  "

# Row 13
Set-TextCell $ws.Range("D13") "29.05"
$ws.Range("E13").Value = "  +6.63%  "

# Row 14
$ws.Range("E14").Value = "  +21.46%  "

# Row 15
$ws.Range("D15").Value = "3.101.46"
$ws.Range("E15").Value = "  +2.14%  "

# Row 16
$ws.Range("D16").Value = "64.907.74"
$ws.Range("E16").Value = "  +3.04%  "

# Row 17
$ws.Range("D17").Value = "2.640.54"
$ws.Range("E17").Value = "  +3.25%  "

# Row 18
$ws.Range("E18").Value = "  +3.26%  "

# Row 19
$ws.Range("E19").Value = "  +3.10%  "

# Row 20
Set-TextCell $ws.Range("D20") "351.31"
$ws.Range("E20").Value = "  +2.38%  "

# Row 21
Set-TextCell $ws.Range("D21") "7.35"
$ws.Range("E21").Value = "  +8.35%  "

# Row 22
$ws.Range("E22").Value = "  +0.20%  "

# Row 23
Set-TextCell $ws.Range("D23") "68.22"
$ws.Range("E23").Value = "  +2.28%  "

# Row 24
$ws.Range("E24").Value = "  +4.96%  "

# Row 25
$ws.Range("E25").Value = "  -2.03%  "

# Row 26
$ws.Range("E26").Value = "  +0.00%  "

# Row 27
$ws.Range("E27").Value = "  +1.94%  "

# Row 28
$ws.Range("E28").Value = "  +1.10%  "

# Rows 29-30: PEPE and Binance-PegBSC-USD swap positions (rank in column A is unchanged)
$ws.Range("B29").Value = "PEPE"
$ws.Range("C29").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
Set-TextCell $ws.Range("D29") "0.0₃0949"
$ws.Range("E29").Value = "  +12.21%  "

$ws.Range("B30").Value = "Binance-PegBSC-USD"
$ws.Range("C30").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
Set-TextCell $ws.Range("D30") "0.999"
$ws.Range("E30").Value = "  -0.25%  "

# Row 31
Set-TextCell $ws.Range("D31") "522.84"
$ws.Range("E31").Value = "  -5.16%  "

# Row 32
$ws.Range("E32").Value = "  +4.44%  "

# Row 33
$ws.Range("E33").Value = "  +1.96%  "

# Row 34
Set-TextCell $ws.Range("D34") "5.55"
$ws.Range("E34").Value = "  +8.06%  "

# Row 35
Set-TextCell $ws.Range("D35") "6.28"
$ws.Range("E35").Value = "  +6.13%  "

# Row 36
$ws.Range("E36").Value = "  +3.50%  "

# Row 37
Set-TextCell $ws.Range("D37") "163.61"
$ws.Range("E37").Value = "  -0.92%  "

# Row 38
Set-TextCell $ws.Range("D38") "20.24"
$ws.Range("E38").Value = "  +4.74%  "

# Row 39
Set-TextCell $ws.Range("D39") "1.98"
$ws.Range("E39").Value = "  +5.80%  "

# Row 40
$ws.Range("E40").Value = "  -0.12%  "

# Row 41
Set-TextCell $ws.Range("D41") "0.999"
$ws.Range("E41").Value = "  +0.03%  "

# Row 42
Set-TextCell $ws.Range("D42") "42.25"
$ws.Range("E42").Value = "  +6.85%  "

# Row 43
Set-TextCell $ws.Range("D43") "164.81"
$ws.Range("E43").Value = "  -0.21%  "

# Row 44
$ws.Range("E44").Value = "  +3.75%  "

# Row 45
$ws.Range("E45").Value = "  +5.11%  "

# Row 46
$ws.Range("E46").Value = "  +2.77%  "

# Row 47
Set-TextCell $ws.Range("D47") "2.19"
$ws.Range("E47").Value = "  +8.59%  "

# Row 48
Set-TextCell $ws.Range("D48") "0.644"
$ws.Range("E48").Value = "  +2.75%  "

# Row 49
$ws.Range("E49").Value = "  +3.34%  "

# Row 50
$ws.Range("E50").Value = "  +2.09%  "

# Row 51
Set-TextCell $ws.Range("D51") "19.41"
$ws.Range("E51").Value = "  +2.93%  "
